$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old table contents (columns B:K, rows 1:2) since the new
# data only occupies A1:A2.
$ws.Range("B1:K2").Clear()

# Set the new values for A1 (header) and A2 (data).
$ws.Range("A1").Value = "Father_name"
$ws.Range("A2").Value = "Vivek"
